# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" positioned right after "总计" (i.e.
#    before the existing "2021-Q2" sheet), carrying the same layout/styles
#    as the other quarterly detail sheets.
# 2. Populate its header + single data row with the 2022-Q3 fund data.
# 3. Update the "总计" (summary) sheet: the new 2022-Q3 row is inserted
#    above the existing 2021-Q2 / 2020-Q4 rows, which shift down by one
#    and have their running index (column A) bumped accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" detail sheet by duplicating the existing
# "2021-Q2" sheet (index 2) - this keeps header styles / column widths /
# borders identical to its sibling quarter sheets. The copy is placed
# immediately before its source, landing at index 2.
# ---------------------------------------------------------------------
$quarterTemplate = $wb.Worksheets.Item(2)
$quarterTemplate.Copy($quarterTemplate)

$newQuarter = $wb.Worksheets.Item(2)
$newQuarter.Name = "2022-Q3"

# the template sheet has 2 data rows; 2022-Q3 only has 1, so drop row 3
$newQuarter.Rows.Item(3).Delete()

# Header row: only the "基金规模" column differs from the template's
# "基金金额" wording.
$newQuarter.Range("D1").Value = "基金规模"

# Data row 2
$newQuarter.Range("B2").Value = "'002020"
$newQuarter.Range("B2").Style = "Normal"
$newQuarter.Range("C2").Value = "国都创新驱动灵活配置混合"
$newQuarter.Range("D2").Value = "'0.12"
$newQuarter.Range("D2").Style = "Normal"
$newQuarter.Range("E2").Value = "'65.45"
$newQuarter.Range("E2").Style = "Normal"
$newQuarter.Range("F2").Value = "'3.15"
$newQuarter.Range("F2").Style = "Normal"
$newQuarter.Range("G2").Value = "'0.0038"
$newQuarter.Range("G2").Style = "Normal"
$newQuarter.Range("H2").Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - shift the two existing rows
# down and write the new 2022-Q3 row in their place.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# row 3 (2020-Q4) -> row 4, bump its running index
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A4").Value = 2

# row 2 (2021-Q2) -> row 3, bump its running index
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))
$summary.Range("A3").Value = 1

# new row 2 holds the 2022-Q3 totals
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# Restore the original "last sheet selected" UI state (the source file
# had its final quarter sheet active/selected rather than the summary).
$wb.Worksheets.Item(4).Activate()
